{"js": "// Change \"Micro I and II (Game Theory)\" -> \"Microeconomics I and II (Game Theory)\"\n// and \"Macro I and II\" -> \"Macroeconomics I and II\" by inserting \"economics\"\n// right after the \"Micro\"/\"Macro\" word in each of the two relevant paragraphs\n// (leaving the unrelated \"International Finance and Macroeconomics\" paragraph\n// untouched).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n\n  if (text === \"Micro I and II (Game Theory)\") {\n    const matches = paragraph.search(\"Micro\", { matchCase: true });\n    matches.load(\"items\");\n    await context.sync();\n    if (matches.items.length > 0) {\n      matches.items[0].insertText(\"economics\", Word.InsertLocation.after);\n    }\n  } else if (text === \"Macro I and II\") {\n    const matches = paragraph.search(\"Macro\", { matchCase: true });\n    matches.load(\"items\");\n    await context.sync();\n    if (matches.items.length > 0) {\n      matches.items[0].insertText(\"economics\", Word.InsertLocation.after);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Change \"Micro I and II (Game Theory)\" -> \"Microeconomics I and II (Game Theory)\"\n# and \"Macro I and II\" -> \"Macroeconomics I and II\" by inserting \"economics\"\n# right after the \"Micro\"/\"Macro\" word in each of the two relevant paragraphs\n# (leaving the unrelated \"International Finance and Macroeconomics\" paragraph\n# untouched).\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $trimmed = $t.TrimEnd([char]13, [char]7)\n\n    if ($trimmed -eq \"Micro I and II (Game Theory)\") {\n        $r = $p.Range\n        $find = $r.Find\n        $find.Text = \"Micro\"\n        $find.MatchCase = $true\n        if ($find.Execute()) {\n            $r.InsertAfter(\"economics\")\n        }\n    }\n    elseif ($trimmed -eq \"Macro I and II\") {\n        $r = $p.Range\n        $find = $r.Find\n        $find.Text = \"Macro\"\n        $find.MatchCase = $true\n        if ($find.Execute()) {\n            $r.InsertAfter(\"economics\")\n        }\n    }\n}\n"}
